$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated quarterly revision to the last existing row (01-01-2021)
$ws.Range("L75").Value = 12116
$ws.Range("M75").Value = 12208

# New quarter row: 01-04-2021
# Enter the period label with a leading apostrophe so it is stored as literal
# text (matching the other period labels in column A) instead of being
# auto-converted to a date serial, then strip the resulting cell formatting
# so the cell keeps the same (default) look as its neighbours.
$ws.Range("A76").Value = "'01-04-2021"
$ws.Range("A76").Style = "Normal"

$ws.Range("B76").Value = 63153
$ws.Range("C76").Value = 6620
$ws.Range("D76").Value = 778
$ws.Range("E76").Value = 55754
$ws.Range("F76").Value = 49871
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 10366
$ws.Range("I76").Value = 6612
$ws.Range("J76").Value = 7659
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 11745
$ws.Range("M76").Value = 11578
$ws.Range("N76").Value = 1912
